# gene_review.xlsx update
# - Apply an AutoFilter on column D ("C") of the data table, showing only
#   rows where the value is FALSE (this hides the TRUE rows + the blank
#   rows below the filtered table, matching the workbook's filtered view).
# - Move the active selection to D25 (single cell).
# - Stamp the sheet with the "OFFICIAL" classification header/footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the AutoFilter: column 4 (D) of the A1:L117 range, keep rows whose
# value is FALSE. Using xlFilterValues (7) with an explicit values array
# produces the <filters><filter val="FALSE"/></filters> representation.
$null = $ws.Range("A1:L117").AutoFilter(4, @("FALSE"), 7)

# Update the selected cell/range to D25.
$null = $ws.Range("D25").Select()

# Add the "OFFICIAL" protective marking to the header and footer.
$ws.PageSetup.CenterHeader = "&`"Calibri`"&12&KFF0000 OFFICIAL&1#`r"
$ws.PageSetup.CenterFooter = "`r&1#&`"Calibri`"&12&KFF0000 OFFICIAL"
